$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 161380
$ws.Range("C4").Value = 152393
$ws.Range("C5").Value = 8987
$ws.Range("C8").Value = 64.63
